$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -1.55324362309513
$ws.Range("D2").Value = 0.1346333785478824

$ws.Range("C3").Value = 0.7774180774052073
$ws.Range("D3").Value = 0.4451920976289536

$ws.Range("C4").Value = 0.2823055712839438
$ws.Range("D4").Value = 0.7803480514046477

$ws.Range("C5").Value = -0.8606525048558071
$ws.Range("D5").Value = 0.3987122780728476

$ws.Range("C6").Value = 1.809929420912749
$ws.Range("D6").Value = 0.08398987302679695
$ws.Range("G6").Value = "No"

$ws.Range("C7").Value = 1.928741580409353
$ws.Range("D7").Value = 0.0667644908603513
$ws.Range("G7").Value = "No"

$ws.Range("C8").Value = 1.122257978133634
$ws.Range("D8").Value = 0.2738522188084636

$ws.Range("C9").Value = -0.4909019660346601
$ws.Range("D9").Value = 0.6283569104521249

$ws.Range("C10").Value = -1.200408582264379
$ws.Range("D10").Value = 0.2427473742636959

$ws.Range("C11").Value = -0.9560022139247716
$ws.Range("D11").Value = 0.3494591377916509
